{"js": "// Auto-generated: replace each old text value with its new value.\n// Each (old, new) pair corresponds to a single <w:t> run that changed\n// in the diff (the date heading + 100 table-cell arithmetic answers).\nconst replacements = [\n  [\"2023-04-20 Thursday\", \"2023-04-21 Friday\"],\n  [\"43-18=25\", \"19+52=71\"],\n  [\"82-36=46\", \"44+8=52\"],\n  [\"76-17=59\", \"44+20=64\"],\n  [\"46+32=78\", \"39-32=7\"],\n  [\"24+26=50\", \"91-62=29\"],\n  [\"70-44=26\", \"12+37=49\"],\n  [\"77-25=52\", \"1+8=9\"],\n  [\"93-76=17\", \"23+76=99\"],\n  [\"95-25=70\", \"90-84=6\"],\n  [\"7+61=68\", \"69-7=62\"],\n  [\"34-0=34\", \"38-0=38\"],\n  [\"33+57=90\", \"54-40=14\"],\n  [\"30-24=6\", \"95-77=18\"],\n  [\"3+57=60\", \"18+53=71\"],\n  [\"76-16=60\", \"80-54=26\"],\n  [\"27-13=14\", \"52+14=66\"],\n  [\"51+13=64\", \"20+49=69\"],\n  [\"35+28=63\", \"58-31=27\"],\n  [\"23+16=39\", \"4+52=56\"],\n  [\"40-32=8\", \"84-32=52\"],\n  [\"62-24=38\", \"33+11=44\"],\n  [\"16-15=1\", \"17+8=25\"],\n  [\"57+35=92\", \"10+68=78\"],\n  [\"58-21=37\", \"78+7=85\"],\n  [\"25-16=9\", \"12+62=74\"],\n  [\"11+50=61\", \"11+40=51\"],\n  [\"25+65=90\", \"61+8=69\"],\n  [\"82-20=62\", \"97-7=90\"],\n  [\"92-4=88\", \"72-32=40\"],\n  [\"37-23=14\", \"77+3=80\"],\n  [\"33+52=85\", \"50+25=75\"],\n  [\"60+19=79\", \"64+13=77\"],\n  [\"67-22=45\", \"2+79=81\"],\n  [\"22-21=1\", \"89-9=80\"],\n  [\"93-14=79\", \"11+37=48\"],\n  [\"47-35=12\", \"10+42=52\"],\n  [\"1+92=93\", \"24+55=79\"],\n  [\"12+23=35\", \"78-39=39\"],\n  [\"53+46=99\", \"26+19=45\"],\n  [\"39+48=87\", \"9+72=81\"],\n  [\"76-15=61\", \"43+13=56\"],\n  [\"75-40=35\", \"29+0=29\"],\n  [\"50-6=44\", \"81-45=36\"],\n  [\"37+14=51\", \"50-38=12\"],\n  [\"33+51=84\", \"74+2=76\"],\n  [\"46-5=41\", \"70-34=36\"],\n  [\"24+52=76\", \"59-47=12\"],\n  [\"85-79=6\", \"95-87=8\"],\n  [\"27+15=42\", \"71-59=12\"],\n  [\"16+2=18\", \"42-11=31\"],\n  [\"54-27=27\", \"83-31=52\"],\n  [\"67-0=67\", \"2+83=85\"],\n  [\"48+48=96\", \"80-29=51\"],\n  [\"25-3=22\", \"22-13=9\"],\n  [\"44+27=71\", \"78+19=97\"],\n  [\"72+9=81\", \"41-26=15\"],\n  [\"0+12=12\", \"27-5=22\"],\n  [\"83-35=48\", \"12-7=5\"],\n  [\"47-5=42\", \"75-65=10\"],\n  [\"30-29=1\", \"10+88=98\"],\n  [\"74+18=92\", \"43+15=58\"],\n  [\"61-10=51\", \"53+23=76\"],\n  [\"84-16=68\", \"88-20=68\"],\n  [\"76-12=64\", \"5+83=88\"],\n  [\"11+70=81\", \"41+46=87\"],\n  [\"38+34=72\", \"73-22=51\"],\n  [\"42+27=69\", \"96-48=48\"],\n  [\"79-42=37\", \"92-19=73\"],\n  [\"88-23=65\", \"11+53=64\"],\n  [\"13+60=73\", \"11+7=18\"],\n  [\"91-30=61\", \"49-45=4\"],\n  [\"51-34=17\", \"45+45=90\"],\n  [\"63-6=57\", \"63-20=43\"],\n  [\"86-0=86\", \"78+4=82\"],\n  [\"71+12=83\", \"75-1=74\"],\n  [\"80-27=53\", \"72-39=33\"],\n  [\"54+25=79\", \"1+82=83\"],\n  [\"79-40=39\", \"32+37=69\"],\n  [\"79+20=99\", \"12-7=5\"],\n  [\"67-4=63\", \"88-62=26\"],\n  [\"34+7=41\", \"30+49=79\"],\n  [\"36+52=88\", \"47+48=95\"],\n  [\"89-61=28\", \"74-47=27\"],\n  [\"84-63=21\", \"11-9=2\"],\n  [\"39-0=39\", \"85-50=35\"],\n  [\"42+3=45\", \"33+66=99\"],\n  [\"56-18=38\", \"30+27=57\"],\n  [\"31+30=61\", \"24+67=91\"],\n  [\"4+82=86\", \"29+26=55\"],\n  [\"84+3=87\", \"12+44=56\"],\n  [\"22+23=45\", \"56+28=84\"],\n  [\"1+94=95\", \"34+18=52\"],\n  [\"86-29=57\", \"35-21=14\"],\n  [\"98-19=79\", \"77-47=30\"],\n  [\"21+73=94\", \"86-69=17\"],\n  [\"48+13=61\", \"75-71=4\"],\n  [\"29+33=62\", \"22+74=96\"],\n  [\"56-5=51\", \"60-6=54\"],\n  [\"0+75=75\", \"0+74=74\"],\n  [\"52+15=67\", \"19+79=98\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Auto-generated: replace each old text value with its new value.\n# Each (old, new) pair corresponds to a single run of text that changed\n# in the diff (the date heading + 100 table-cell arithmetic answers).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@('2023-04-20 Thursday', '2023-04-21 Friday')\n    ,@('43-18=25', '19+52=71')\n    ,@('82-36=46', '44+8=52')\n    ,@('76-17=59', '44+20=64')\n    ,@('46+32=78', '39-32=7')\n    ,@('24+26=50', '91-62=29')\n    ,@('70-44=26', '12+37=49')\n    ,@('77-25=52', '1+8=9')\n    ,@('93-76=17', '23+76=99')\n    ,@('95-25=70', '90-84=6')\n    ,@('7+61=68', '69-7=62')\n    ,@('34-0=34', '38-0=38')\n    ,@('33+57=90', '54-40=14')\n    ,@('30-24=6', '95-77=18')\n    ,@('3+57=60', '18+53=71')\n    ,@('76-16=60', '80-54=26')\n    ,@('27-13=14', '52+14=66')\n    ,@('51+13=64', '20+49=69')\n    ,@('35+28=63', '58-31=27')\n    ,@('23+16=39', '4+52=56')\n    ,@('40-32=8', '84-32=52')\n    ,@('62-24=38', '33+11=44')\n    ,@('16-15=1', '17+8=25')\n    ,@('57+35=92', '10+68=78')\n    ,@('58-21=37', '78+7=85')\n    ,@('25-16=9', '12+62=74')\n    ,@('11+50=61', '11+40=51')\n    ,@('25+65=90', '61+8=69')\n    ,@('82-20=62', '97-7=90')\n    ,@('92-4=88', '72-32=40')\n    ,@('37-23=14', '77+3=80')\n    ,@('33+52=85', '50+25=75')\n    ,@('60+19=79', '64+13=77')\n    ,@('67-22=45', '2+79=81')\n    ,@('22-21=1', '89-9=80')\n    ,@('93-14=79', '11+37=48')\n    ,@('47-35=12', '10+42=52')\n    ,@('1+92=93', '24+55=79')\n    ,@('12+23=35', '78-39=39')\n    ,@('53+46=99', '26+19=45')\n    ,@('39+48=87', '9+72=81')\n    ,@('76-15=61', '43+13=56')\n    ,@('75-40=35', '29+0=29')\n    ,@('50-6=44', '81-45=36')\n    ,@('37+14=51', '50-38=12')\n    ,@('33+51=84', '74+2=76')\n    ,@('46-5=41', '70-34=36')\n    ,@('24+52=76', '59-47=12')\n    ,@('85-79=6', '95-87=8')\n    ,@('27+15=42', '71-59=12')\n    ,@('16+2=18', '42-11=31')\n    ,@('54-27=27', '83-31=52')\n    ,@('67-0=67', '2+83=85')\n    ,@('48+48=96', '80-29=51')\n    ,@('25-3=22', '22-13=9')\n    ,@('44+27=71', '78+19=97')\n    ,@('72+9=81', '41-26=15')\n    ,@('0+12=12', '27-5=22')\n    ,@('83-35=48', '12-7=5')\n    ,@('47-5=42', '75-65=10')\n    ,@('30-29=1', '10+88=98')\n    ,@('74+18=92', '43+15=58')\n    ,@('61-10=51', '53+23=76')\n    ,@('84-16=68', '88-20=68')\n    ,@('76-12=64', '5+83=88')\n    ,@('11+70=81', '41+46=87')\n    ,@('38+34=72', '73-22=51')\n    ,@('42+27=69', '96-48=48')\n    ,@('79-42=37', '92-19=73')\n    ,@('88-23=65', '11+53=64')\n    ,@('13+60=73', '11+7=18')\n    ,@('91-30=61', '49-45=4')\n    ,@('51-34=17', '45+45=90')\n    ,@('63-6=57', '63-20=43')\n    ,@('86-0=86', '78+4=82')\n    ,@('71+12=83', '75-1=74')\n    ,@('80-27=53', '72-39=33')\n    ,@('54+25=79', '1+82=83')\n    ,@('79-40=39', '32+37=69')\n    ,@('79+20=99', '12-7=5')\n    ,@('67-4=63', '88-62=26')\n    ,@('34+7=41', '30+49=79')\n    ,@('36+52=88', '47+48=95')\n    ,@('89-61=28', '74-47=27')\n    ,@('84-63=21', '11-9=2')\n    ,@('39-0=39', '85-50=35')\n    ,@('42+3=45', '33+66=99')\n    ,@('56-18=38', '30+27=57')\n    ,@('31+30=61', '24+67=91')\n    ,@('4+82=86', '29+26=55')\n    ,@('84+3=87', '12+44=56')\n    ,@('22+23=45', '56+28=84')\n    ,@('1+94=95', '34+18=52')\n    ,@('86-29=57', '35-21=14')\n    ,@('98-19=79', '77-47=30')\n    ,@('21+73=94', '86-69=17')\n    ,@('48+13=61', '75-71=4')\n    ,@('29+33=62', '22+74=96')\n    ,@('56-5=51', '60-6=54')\n    ,@('0+75=75', '0+74=74')\n    ,@('52+15=67', '19+79=98')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
